# Update the "dSF" column (F) values for multiple rows, as per the
# repull/recalculation of data described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 2
    4  = 4
    5  = 3
    6  = 3
    7  = -4
    8  = -2
    9  = -1
    11 = 2
    12 = -3
    13 = -3
    14 = 2
    15 = -4
    16 = -1
    17 = -1
    18 = 1
    19 = 3
    20 = -4
    21 = 3
    24 = -2
    25 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
